$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Expand the header row (row 1) on the "Centers" sheet ---
# Copy the existing header style (bold + border, style index 1) from E1 onto
# the newly added header cells F1:K1 so they visually match B1:E1.
$ws1.Range("E1").Copy() | Out-Null
$ws1.Range("F1:K1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# New header labels (columns F through K)
$ws1.Range("F1").Value = "Facility City Name"
$ws1.Range("G1").Value = "Facility City Name"
$ws1.Range("H1").Value = "Facility Postal Code"
$ws1.Range("I1").Value = "Deliv Center Capac"
$ws1.Range("J1").Value = "Latitude"
$ws1.Range("K1").Value = "Longitude"

# Re-label the pre-existing header columns C1, D1, E1 (B1 "Facility Name" stays)
$ws1.Range("C1").Value = "Center Name"
$ws1.Range("D1").Value = "Center Num"
$ws1.Range("E1").Value = "Facility Address Line 1"

# --- Row 2 (MYKAWA) ---
$ws1.Range("C2").Value = "MYKAWA EAM"
$ws1.Range("D2").Value = 7723
$ws1.Range("E2").Value = "7110 MYKAWA ROAD"
$ws1.Range("F2").Value = "HOUSTON"
$ws1.Range("G2").Value = "HOUSTON"
$ws1.Range("H2").Value = 77033
$ws1.Range("I2").Value = 999999999
$ws1.Range("J2").Value = 29.67578534220857
$ws1.Range("K2").Value = -95.32125610590822

# --- Row 3 (STAFFORD) ---
$ws1.Range("C3").Value = "STAFFORD-KATY"
$ws1.Range("D3").Value = 7741
$ws1.Range("E3").Value = "13922 STAFFORD ROAD"
$ws1.Range("F3").Value = "STAFFORD"
$ws1.Range("G3").Value = "STAFFORD"
$ws1.Range("H3").Value = 77477
$ws1.Range("I3").Value = 999999999
$ws1.Range("J3").Value = 29.61935353306665
$ws1.Range("K3").Value = -95.55583702616946

# --- Row 4 (HOUSTON (CANINO) HUB) ---
$ws1.Range("C4").Value = "SWEETWATER-CENTRAL"
$ws1.Range("D4").Value = 7707
$ws1.Range("E4").Value = "8330 SWEETWATER LANE"
$ws1.Range("F4").Value = "HOUSTON"
$ws1.Range("G4").Value = "HOUSTON"
$ws1.Range("H4").Value = 77037
$ws1.Range("I4").Value = 999999999
$ws1.Range("J4").Value = 29.88062033368866
$ws1.Range("K4").Value = -95.40848234626999
